$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the crypto price refresh diff.
# "D" column (Price) cells are forced to Text so numeric-looking
# strings (e.g. "306.02") are not reinterpreted as numbers by Excel,
# matching the original inlineStr text cells, then formats are reset
# so no stray style index is left on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.452.16"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.810.67"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "306.02"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4501"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3582"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.40"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07055"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8911"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07798"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.36"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.838.89"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.274"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.306"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "84.79"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008529"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.006"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.476.33"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.19"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.963"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.064.33"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.77%  "
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.953"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.57"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.054"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "112.18"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08682"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.114"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.64%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7385"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.53%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.754"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +8.88%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.111"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.072"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01923"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05124"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.891"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5088"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.761"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.64%  "
$ws.Range("E44").Value = "  -3.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.044"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4675"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.006"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.02"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "99.68"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.80%  "
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05992"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.15%  "
